# Amazon locators workbook update:
#  - remove headphones second-product locator, replace with MacBook search result locator
#  - add new MacBook-order related locators / verification entities
#  - widen a couple of columns to fit the new (longer) locator strings
#  - update the saved selection / active cell

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Locator_ids
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Locator_ids")

# Row 15 ("Amazon_Search_Second_Product_xpath") now points at the MacBook
# search result instead of the generic second-product image.
$ws1.Range("B15").Value = "xpath=(//div[contains(@data-cel-widget,'search_result_3')]//span[contains(text(),'MacBook')])"

# New rows appended for the MacBook add-to-cart / checkout flow.
$ws1.Range("A25").Value = "Amazon_Product_MacBook_Titles_Added_In_Cart_xpath"
$ws1.Range("B25").Value = "xpath=(//span[contains(@class,'sc-product-title')][contains(text(),'MacBook')])"

$ws1.Range("A26").Value = "Amazon_Macbook_Order_Close_Button_xpath"
$ws1.Range("B26").Value = "xpath=(//button[contains(@class,'button-close')]/i)"

$ws1.Range("A27").Value = "Amazon_Macbook_Order_Continue_Button_xpath"
$ws1.Range("B27").Value = "xpath=//span[contains(text(),'Continue')]/..//input[contains(@class,'a-button-input')]"
$ws1.Range("H27").Value = "xpath=(//span[contains(@class,'button')][contains(text(),'Continue')])/..//input"

# Column B needs to be wider to fit the new, longer locator strings.
$ws1.Columns.Item(2).ColumnWidth = 88.33333333333334

$ws1.Activate()
$ws1.Range("B27").Select()

# ---------------------------------------------------------------------
# Sheet 2: Verification_Entities
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Verification_Entities")

$ws2.Range("A4").Value = "Verify_Amazon_Added_To_Cart_Text_Message_Macbook_xpath"
$ws2.Range("B4").Value = "xpath=(//h1[contains(text(),'Added to Cart')])"

$ws2.Range("B5").Value = "xpath=(//span[contains(@class,'button')][contains(text(),'Continue')])"

# Column A needs to be wider too.
$ws2.Columns.Item(1).ColumnWidth = 68.66666666666666

$ws2.Activate()
$ws2.Range("B5").Select()

# Leave Sheet 1 as the active / tab-selected sheet, matching the source file.
$ws1.Activate()
